# Fill in the "owner_id" (column C) and "manager_id" (column D) values
# for the first few data rows (institutions and folders source codes).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "ofs"
$ws.Range("D2").Value = "ofs-tourisme"

$ws.Range("C3").Value = "ofs"
$ws.Range("D3").Value = "ofs-div-pop"

$ws.Range("C4").Value = "ofs"
$ws.Range("D4").Value = "ofs-travail"

# Reflect the final selection/active cell recorded in the saved file.
$ws.Range("E15").Select()
